# Apiary certificate: remove the signatory's name from the footer and
# bump the "dated at" date shown next to it.
#
# wdReplaceAll = 2, wdFindContinue = 1 (Wrap parameter)
$wdReplaceAll = 2
$wdFindContinue = 1

$d = $word.ActiveDocument

# The signature block ("ABBOTSFORD ... <name> , <date> ... Provincial
# Apiculturist") lives in the primary/default footer (Footers(1) ==
# wdHeaderFooterPrimary) of the document's only section.
$footerRange = $d.Sections(1).Footers(1).Range

# Drop "Paul van Westendorp," entirely - the trailing tab that used to
# precede the name is left in place, the paragraph now ends right after it.
$footerRange.Find.Execute("Paul van Westendorp,", $true, $false, $false, `
    $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceAll)

# Re-scope and update the cached result of the "DATE \@ yyyy-MM-dd" field
# next to the (now-removed) signature so the certificate shows the new
# issue date.
$footerRange = $d.Sections(1).Footers(1).Range
$footerRange.Find.Execute("2023-03-01", $true, $false, $false, `
    $false, $false, $true, $wdFindContinue, $false, "2025-10-01", $wdReplaceAll)
